$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Search"
$ws.Range("B3").Value = "Search + SB"
$ws.Range("B4").Value = "GS + SB"
$ws.Range("B5").Value = "Scopus + SB"
$ws.Range("B6").Value = "Scopus + BW // FW"
$ws.Range("B7").Value = "Scopus + BW + FW"
$ws.Range("B8").Value = "Scopus + FW + BW"
